$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update in place ---
$ws.Range("A2").Value = "lala morale"
$ws.Range("B2").Value = "'113564"
$ws.Range("C2").Value = "'114321654687987654543213"
$ws.Range("D2").Value = "'"
$ws.Range("E2").Value = "'23132"
$ws.Range("F2").Value = "Direction régionale"
$ws.Range("G2").Value = "001/TEST DR/AV"
$ws.Range("H2").Value = "mensuelle"
$ws.Range("I2").Value = 48000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 48000

# --- Row 3: update in place ---
$ws.Range("A3").Value = "YASSINE TYEST"
$ws.Range("B3").Value = "BB125874"
$ws.Range("C3").Value = "'115649679785432432321321"
$ws.Range("D3").Value = "'"
$ws.Range("E3").Value = "tesqt"
$ws.Range("F3").Value = "Direction régionale"
$ws.Range("G3").Value = "001/TEST DR/AV"
$ws.Range("H3").Value = "mensuelle"
$ws.Range("I3").Value = 72000
$ws.Range("J3").Value = 5400
$ws.Range("K3").Value = 66600

# --- Remove old rows 4-7 entirely (they are no longer part of the table) ---
$ws.Range("A4:K7").EntireRow.Delete()
